$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (Pollock): add Landings (B) and Target TAC (D) ---
$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 8832

$ws.Range("D4").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 55200

# --- Row 15: rename "Redfish - 500s" -> "Redfish", restyle to match other
#     species headers (bold Arial 9, left/top aligned box), add Landings
#     and Target TAC ---
$ws.Range("A17").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Redfish"
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Font.Size = 9
$ws.Range("A15").Font.Name = "Arial"
$ws.Range("A15").HorizontalAlignment = -4131
$ws.Range("A15").VerticalAlignment = -4160

$ws.Range("B6").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 677

$ws.Range("D4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 67700

# --- Row 17 (White Hake): add Target TAC (D) ---
$ws.Range("D4").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 53890.909090909088

# --- Row 19 (American Plaice): add Target TAC (D) ---
$ws.Range("D4").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 58358.333333333336

# --- Selection moved ---
$ws.Range("H22").Select()

$excel.CutCopyMode = 0
